# Auto-generated edit script applying the Tonberry_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 90.72727
$ws.Range("I33").Value = 70
$ws.Range("J33").Value = 115.6
$ws.Range("K33").Value = 70
$ws.Range("L33").Value = 115.6
$ws.Range("M33").Value = 159
$ws.Range("N33").Value = -573.6
$ws.Range("H76").Value = 3154.6
$ws.Range("I76").Value = 2924.6667
$ws.Range("K76").Value = 2924.6667
$ws.Range("M76").Value = -2609.6667
$ws.Range("H79").Value = 3154.6
$ws.Range("I79").Value = 2924.6667
$ws.Range("K79").Value = 2924.6667
$ws.Range("M79").Value = -1832.6667
$ws.Range("H113").Value = 44786.43
$ws.Range("I113").Value = 75751.25
$ws.Range("K113").Value = 75751.25
$ws.Range("M113").Value = -72497.25
$ws.Range("H123").Value = 37998.75
$ws.Range("J123").Value = 37998.75
$ws.Range("L123").Value = 37998.75
$ws.Range("N123").Value = -47798.75
$ws.Range("H131").Value = 2060.5
$ws.Range("I131").Value = 736.8461
$ws.Range("J131").Value = 5502
$ws.Range("K131").Value = 2210.5383
$ws.Range("L131").Value = 16506
$ws.Range("M131").Value = 2829.4617
$ws.Range("N131").Value = -26586
$ws.Range("H132").Value = 1387.6031
$ws.Range("I132").Value = 1204.2115
$ws.Range("K132").Value = 3612.6345
$ws.Range("M132").Value = -1082.6345
$ws.Range("H135").Value = 295.13513
$ws.Range("I135").Value = 295.13513
$ws.Range("K135").Value = 2656.21617
$ws.Range("M135").Value = -121.2161700000001
$ws.Range("H137").Value = 29688.457
$ws.Range("I137").Value = 823.2917
$ws.Range("J137").Value = 92667
$ws.Range("K137").Value = 2469.8751
$ws.Range("L137").Value = 278001
$ws.Range("M137").Value = 80.1248999999998
$ws.Range("N137").Value = -283101
$ws.Range("H138").Value = 1877.3827
$ws.Range("I138").Value = 1575.6666
$ws.Range("J138").Value = 2480.8147
$ws.Range("K138").Value = 4726.9998
$ws.Range("L138").Value = 7442.4441
$ws.Range("M138").Value = 413.0002000000004
$ws.Range("N138").Value = -17722.4441

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3515.6753
$ws.Range("I32").Value = 3287.0303
$ws.Range("J32").Value = 4887.5454
$ws.Range("K32").Value = 3287.0303
$ws.Range("L32").Value = 4887.5454
$ws.Range("M32").Value = -3000.0303
$ws.Range("N32").Value = -5461.5454
$ws.Range("H61").Value = 2413.75
$ws.Range("I61").Value = 1741.6666
$ws.Range("K61").Value = 1741.6666
$ws.Range("M61").Value = -1529.6666
$ws.Range("H63").Value = 8000
$ws.Range("I63").Value = 8000
$ws.Range("K63").Value = 8000
$ws.Range("M63").Value = -7314
$ws.Range("H66").Value = 8000
$ws.Range("I66").Value = 8000
$ws.Range("K66").Value = 40000
$ws.Range("M66").Value = -36568
$ws.Range("H74").Value = 1761.4
$ws.Range("I74").Value = 1413.5555
$ws.Range("K74").Value = 1413.5555
$ws.Range("M74").Value = -539.5554999999999
$ws.Range("H77").Value = 1761.4
$ws.Range("I77").Value = 1413.5555
$ws.Range("K77").Value = 7067.7775
$ws.Range("M77").Value = -2699.7775
$ws.Range("H132").Value = 1646.4219
$ws.Range("I132").Value = 1129.3191
$ws.Range("J132").Value = 3076.0588
$ws.Range("K132").Value = 3387.9573
$ws.Range("L132").Value = 9228.1764
$ws.Range("M132").Value = -857.9573
$ws.Range("N132").Value = -14288.1764
$ws.Range("H136").Value = 2413.75
$ws.Range("I136").Value = 1741.6666
$ws.Range("K136").Value = 5224.9998
$ws.Range("M136").Value = -2674.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1814.7428
$ws.Range("I20").Value = 1436.96
$ws.Range("J20").Value = 2759.2
$ws.Range("K20").Value = 1436.96
$ws.Range("L20").Value = 2759.2
$ws.Range("M20").Value = -1189.96
$ws.Range("N20").Value = -3253.2
$ws.Range("H35").Value = 35000
$ws.Range("J35").Value = 35000
$ws.Range("L35").Value = 35000
$ws.Range("N35").Value = -35620
$ws.Range("H94").Value = 758.9091
$ws.Range("J94").Value = 278.57144
$ws.Range("L94").Value = 278.57144
$ws.Range("N94").Value = -1180.57144
$ws.Range("H105").Value = 2149.64
$ws.Range("I105").Value = 2106.4092
$ws.Range("J105").Value = 2466.6667
$ws.Range("K105").Value = 2106.4092
$ws.Range("L105").Value = 2466.6667
$ws.Range("M105").Value = -359.4092000000001
$ws.Range("N105").Value = -5960.6667
$ws.Range("H134").Value = 4607.0884
$ws.Range("I134").Value = 4783.4136
$ws.Range("K134").Value = 14350.2408
$ws.Range("M134").Value = -11815.2408

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2127.9412
$ws.Range("I31").Value = 1458.6111
$ws.Range("K31").Value = 1458.6111
$ws.Range("M31").Value = -1163.6111
$ws.Range("H34").Value = 2127.9412
$ws.Range("I34").Value = 1458.6111
$ws.Range("K34").Value = 1458.6111
$ws.Range("M34").Value = -1256.6111
$ws.Range("H58").Value = 888253.2
$ws.Range("I58").Value = 1115457.6
$ws.Range("K58").Value = 1115457.6
$ws.Range("M58").Value = -1115254.6
$ws.Range("H93").Value = 13366.5
$ws.Range("I93").Value = 10549.75
$ws.Range("J93").Value = 19000
$ws.Range("K93").Value = 10549.75
$ws.Range("L93").Value = 19000
$ws.Range("M93").Value = -8677.75
$ws.Range("N93").Value = -22744
$ws.Range("H107").Value = 798.5
$ws.Range("I107").Value = 520.3077
$ws.Range("K107").Value = 520.3077
$ws.Range("M107").Value = 1399.6923
$ws.Range("H132").Value = 1952
$ws.Range("I132").Value = 1344.2858
$ws.Range("J132").Value = 3370
$ws.Range("K132").Value = 4032.8574
$ws.Range("L132").Value = 10110
$ws.Range("M132").Value = -1502.8574
$ws.Range("N132").Value = -15170
$ws.Range("H134").Value = 1974.05
$ws.Range("I134").Value = 1804.5454
$ws.Range("J134").Value = 2773.1428
$ws.Range("K134").Value = 5413.6362
$ws.Range("L134").Value = 8319.428400000001
$ws.Range("M134").Value = -2878.6362
$ws.Range("N134").Value = -13389.4284
$ws.Range("H136").Value = 888253.2
$ws.Range("I136").Value = 1115457.6
$ws.Range("K136").Value = 3346372.8
$ws.Range("M136").Value = -3343822.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 374834.9
$ws.Range("I4").Value = 156
$ws.Range("K4").Value = 468
$ws.Range("M4").Value = -356
$ws.Range("H131").Value = 35772840
$ws.Range("J131").Value = 90773.89
$ws.Range("L131").Value = 272321.67
$ws.Range("N131").Value = -282401.67
$ws.Range("H141").Value = 2922.0625
$ws.Range("J141").Value = 3166.6667
$ws.Range("L141").Value = 9500.000100000001
$ws.Range("N141").Value = -19860.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 43003
$ws.Range("I22").Value = 5997
$ws.Range("K22").Value = 5997
$ws.Range("M22").Value = -5468
$ws.Range("H70").Value = 3044.125
$ws.Range("I70").Value = 2951.2354
$ws.Range("K70").Value = 2951.2354
$ws.Range("M70").Value = -2681.2354
$ws.Range("H73").Value = 3044.125
$ws.Range("I73").Value = 2951.2354
$ws.Range("K73").Value = 2951.2354
$ws.Range("M73").Value = -2015.2354
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H104").Value = 49997
$ws.Range("J104").Value = 49997
$ws.Range("L104").Value = 49997
$ws.Range("N104").Value = -56985
$ws.Range("H122").Value = 1483.7
$ws.Range("I122").Value = 1383.625
$ws.Range("J122").Value = 1884
$ws.Range("K122").Value = 4150.875
$ws.Range("L122").Value = 5652
$ws.Range("M122").Value = -1700.875
$ws.Range("N122").Value = -10552
$ws.Range("H132").Value = 701204.7
$ws.Range("I132").Value = 1041003.25
$ws.Range("K132").Value = 3123009.75
$ws.Range("M132").Value = -3120479.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 18000
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H135").Value = 32719
$ws.Range("J135").Value = 32719
$ws.Range("L135").Value = 32719
$ws.Range("N135").Value = -42859
$ws.Range("H136").Value = 1971.758
$ws.Range("I136").Value = 1381.7255
$ws.Range("K136").Value = 4145.1765
$ws.Range("M136").Value = -1595.1765

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 49999
$ws.Range("J105").Value = 49999
$ws.Range("L105").Value = 49999
$ws.Range("N105").Value = -56987
$ws.Range("H132").Value = 1626.659
$ws.Range("I132").Value = 1101.1333
$ws.Range("K132").Value = 3303.3999
$ws.Range("M132").Value = -773.3998999999999
$ws.Range("H136").Value = 9109146
$ws.Range("I136").Value = 12628062
$ws.Range("J136").Value = 1362.8823
$ws.Range("K136").Value = 37884186
$ws.Range("L136").Value = 4088.6469
$ws.Range("M136").Value = -37881636
$ws.Range("N136").Value = -9188.6469
